# Saldo_guide.xlsx update: refresh the daily balance export from
# 2024-06-06 (IClientBalance-20240606-100000-) to 2024-06-07
# (IClientBalance-20240607-085141-).
#
# - Every row's "Dt. Referencia" (column G) moves from 45449 to 45450
#   (2024-06-06 -> 2024-06-07).
# - A subset of accounts picked up a projected movement ("Vl. Projetado",
#   column E) that shifts their "Vl. Total" (column H = D + E).
# - One account's "Saldo Previsto" (column D) itself was corrected, with
#   "Vl. Total" following it (E stays 0).
# - The sheet (tab) name is renamed to match the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to the new export timestamp.
$ws.Name = "IClientBalance-20240607-085141-"

# All data rows (2-257): advance the reference date by one day.
$ws.Range("G2:G257").Value = 45450

# Rows whose projected value (E) moved away from 0, with Vl. Total (H)
# recomputed as D + E.
$projectedMoves = @{
    5   = @{ E = -11235.26;            H = -10523.71 }
    8   = @{ E = -4544.7700000000004;  H = -3789.3 }
    15  = @{ E = -14220.45;            H = -12540.69 }
    17  = @{ E = -5020.72;             H = -4767.72 }
    42  = @{ E = -6692.4;              H = -5702.56 }
    57  = @{ E = -2036.74;             H = -351.76 }
    59  = @{ E = -8605.2800000000007;  H = -8155.01 }
    98  = @{ E = -8840.42;             H = -8198.2199999999993 }
    103 = @{ E = -24126.61;            H = -23333.66 }
    107 = @{ E = -28187.41;            H = -27132.82 }
    131 = @{ E = -3588.32;             H = -3138.34 }
    141 = @{ E = -31758.65;            H = -31758.1 }
    155 = @{ E = -604.16999999999996;  H = 250.95 }
    168 = @{ E = -1559.1;              H = -1223.1199999999999 }
    226 = @{ E = -7885.49;             H = -7705.99 }
    240 = @{ E = -8124.9;              H = -7651.98 }
}

foreach ($row in $projectedMoves.Keys) {
    $vals = $projectedMoves[$row]
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 8).Value = $vals.H
}

# Row 39: Saldo Previsto (D) corrected directly; Vl. Projetado (E) stays 0,
# so Vl. Total (H) just follows the new D.
$ws.Cells.Item(39, 4).Value = 94.6
$ws.Cells.Item(39, 8).Value = 94.6
